# Auto-generated script applying scraped value updates to Sargatanas_Profits sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H87").Value = 0
$ws.Range("J87").Value = 0
$ws.Range("L87").Value = 0
$ws.Range("N87").Value = $null
$ws.Range("H90").Value = 0
$ws.Range("J90").Value = 0
$ws.Range("L90").Value = 0
$ws.Range("N90").Value = $null
$ws.Range("H116").Value = 19235708
$ws.Range("I116").Value = 62501450
$ws.Range("J116").Value = 6490.5557
$ws.Range("K116").Value = 62501450
$ws.Range("L116").Value = 6490.5557
$ws.Range("M116").Value = -62498008
$ws.Range("N116").Value = -13374.5557
$ws.Range("H131").Value = 1674.4286
$ws.Range("I131").Value = 1203.6666
$ws.Range("K131").Value = 3610.9998
$ws.Range("M131").Value = 1429.0002
$ws.Range("H132").Value = 1510.1904
$ws.Range("I132").Value = 1510.1904
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 4530.5712
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -2000.5712
$ws.Range("N132").Value = $null
$ws.Range("H135").Value = 417297.72
$ws.Range("I135").Value = 526750.9
$ws.Range("K135").Value = 4740758.100000001
$ws.Range("M135").Value = -4738223.100000001
$ws.Range("H137").Value = 2546.423
$ws.Range("J137").Value = 2648.25
$ws.Range("L137").Value = 7944.75
$ws.Range("N137").Value = -13044.75
$ws.Range("H138").Value = 4716.549
$ws.Range("I138").Value = 881.4643
$ws.Range("J138").Value = 9385.348
$ws.Range("K138").Value = 2644.3929
$ws.Range("L138").Value = 28156.044
$ws.Range("M138").Value = 2495.6071
$ws.Range("N138").Value = -38436.044
$ws.Range("H141").Value = 1725.5625
$ws.Range("I141").Value = 1710.6666
$ws.Range("J141").Value = 1770.25
$ws.Range("K141").Value = 5131.9998
$ws.Range("L141").Value = 5310.75
$ws.Range("M141").Value = 48.0002000000004
$ws.Range("N141").Value = -15670.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5005973.5
$ws.Range("I32").Value = 5687744
$ws.Range("J32").Value = 6322.6665
$ws.Range("K32").Value = 5687744
$ws.Range("L32").Value = 6322.6665
$ws.Range("M32").Value = -5687457
$ws.Range("N32").Value = -6896.6665
$ws.Range("H45").Value = 2210.1667
$ws.Range("I45").Value = 1948.5
$ws.Range("K45").Value = 1948.5
$ws.Range("M45").Value = -1571.5
$ws.Range("H46").Value = 4999
$ws.Range("J46").Value = 4999
$ws.Range("L46").Value = 4999
$ws.Range("N46").Value = -5637
$ws.Range("H61").Value = 4350.0557
$ws.Range("I61").Value = 3231.3845
$ws.Range("J61").Value = 7258.6
$ws.Range("K61").Value = 3231.3845
$ws.Range("L61").Value = 7258.6
$ws.Range("M61").Value = -3019.3845
$ws.Range("N61").Value = -7682.6
$ws.Range("H74").Value = 46434.027
$ws.Range("I74").Value = 56427.484
$ws.Range("K74").Value = 56427.484
$ws.Range("M74").Value = -55553.484
$ws.Range("H77").Value = 46434.027
$ws.Range("I77").Value = 56427.484
$ws.Range("K77").Value = 282137.42
$ws.Range("M77").Value = -277769.42
$ws.Range("H136").Value = 4350.0557
$ws.Range("I136").Value = 3231.3845
$ws.Range("J136").Value = 7258.6
$ws.Range("K136").Value = 9694.1535
$ws.Range("L136").Value = 21775.8
$ws.Range("M136").Value = -7144.1535
$ws.Range("N136").Value = -26875.8

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H68").Value = 25000
$ws.Range("I68").Value = 25000
$ws.Range("K68").Value = 25000
$ws.Range("M68").Value = -24189
$ws.Range("H71").Value = 25000
$ws.Range("I71").Value = 25000
$ws.Range("K71").Value = 75000
$ws.Range("M71").Value = -70944
$ws.Range("H94").Value = 1979.7142
$ws.Range("I94").Value = 732.8889
$ws.Range("K94").Value = 732.8889
$ws.Range("M94").Value = -281.8889
$ws.Range("H134").Value = 3703.0256
$ws.Range("I134").Value = 1986.5781
$ws.Range("K134").Value = 5959.7343
$ws.Range("M134").Value = -3424.7343

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6614.971
$ws.Range("I31").Value = 2879.8718
$ws.Range("K31").Value = 2879.8718
$ws.Range("M31").Value = -2584.8718
$ws.Range("H34").Value = 6614.971
$ws.Range("I34").Value = 2879.8718
$ws.Range("K34").Value = 2879.8718
$ws.Range("M34").Value = -2677.8718
$ws.Range("H132").Value = 8824.25
$ws.Range("I132").Value = 2378.2
$ws.Range("K132").Value = 7134.599999999999
$ws.Range("M132").Value = -4604.599999999999
$ws.Range("H134").Value = 8310.021000000001
$ws.Range("I134").Value = 7880.25
$ws.Range("J134").Value = 8778.862999999999
$ws.Range("K134").Value = 23640.75
$ws.Range("L134").Value = 26336.589
$ws.Range("M134").Value = -21105.75
$ws.Range("N134").Value = -31406.589

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 13985883
$ws.Range("I4").Value = 15979931
$ws.Range("J4").Value = 27550
$ws.Range("K4").Value = 47939793
$ws.Range("L4").Value = 82650
$ws.Range("M4").Value = -47939681
$ws.Range("N4").Value = -82874
$ws.Range("H23").Value = 210.71428
$ws.Range("I23").Value = 172.5
$ws.Range("J23").Value = 261.66666
$ws.Range("K23").Value = 517.5
$ws.Range("L23").Value = 784.9999799999999
$ws.Range("M23").Value = -282.5
$ws.Range("N23").Value = -1254.99998
$ws.Range("H68").Value = 4310.2856
$ws.Range("I68").Value = 1329.8
$ws.Range("K68").Value = 3989.4
$ws.Range("M68").Value = -3178.4
$ws.Range("H71").Value = 4310.2856
$ws.Range("I71").Value = 1329.8
$ws.Range("K71").Value = 11968.2
$ws.Range("M71").Value = -7912.199999999999
$ws.Range("H131").Value = 1520
$ws.Range("J131").Value = 2543.1667
$ws.Range("L131").Value = 7629.500100000001
$ws.Range("N131").Value = -17709.5001
$ws.Range("H132").Value = 15391.223
$ws.Range("I132").Value = 11724.6
$ws.Range("J132").Value = 19974.5
$ws.Range("K132").Value = 105521.4
$ws.Range("L132").Value = 179770.5
$ws.Range("M132").Value = -102991.4
$ws.Range("N132").Value = -184830.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 2272.8708
$ws.Range("I97").Value = 2050
$ws.Range("J97").Value = 2581.4614
$ws.Range("K97").Value = 2050
$ws.Range("L97").Value = 2581.4614
$ws.Range("M97").Value = -1554
$ws.Range("N97").Value = -3573.4614
$ws.Range("H122").Value = 58070.42
$ws.Range("I122").Value = 170065.83
$ws.Range("J122").Value = 6380.231
$ws.Range("K122").Value = 510197.49
$ws.Range("L122").Value = 19140.693
$ws.Range("M122").Value = -507747.49
$ws.Range("N122").Value = -24040.693
$ws.Range("H132").Value = 2506.543
$ws.Range("I132").Value = 2469.6562
$ws.Range("J132").Value = 2900
$ws.Range("K132").Value = 7408.9686
$ws.Range("L132").Value = 8700
$ws.Range("M132").Value = -4878.9686
$ws.Range("N132").Value = -13760

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H17").Value = 500
$ws.Range("I17").Value = 500
$ws.Range("K17").Value = 500
$ws.Range("M17").Value = -330
$ws.Range("H40").Value = 4899.6665
$ws.Range("I40").Value = 4364.35
$ws.Range("K40").Value = 4364.35
$ws.Range("M40").Value = -4228.35
$ws.Range("H55").Value = 52631908
$ws.Range("I55").Value = 250000060
$ws.Range("J55").Value = 401.06668
$ws.Range("K55").Value = 250000060
$ws.Range("L55").Value = 401.06668
$ws.Range("M55").Value = -249999887
$ws.Range("N55").Value = -747.06668
$ws.Range("H93").Value = 6353.0835
$ws.Range("I93").Value = 5125.125
$ws.Range("K93").Value = 5125.125
$ws.Range("M93").Value = -3877.125
$ws.Range("H132").Value = 17865672
$ws.Range("I132").Value = 33338366
$ws.Range("K132").Value = 100015098
$ws.Range("M132").Value = -100012568
$ws.Range("H133").Value = 150000
$ws.Range("J133").Value = 150000
$ws.Range("L133").Value = 150000
$ws.Range("N133").Value = -155060
$ws.Range("H136").Value = 7666.727
$ws.Range("I136").Value = 4153.28
$ws.Range("K136").Value = 12459.84
$ws.Range("M136").Value = -9909.84

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 1570.6
$ws.Range("I107").Value = 1701
$ws.Range("J107").Value = 1375
$ws.Range("K107").Value = 5103
$ws.Range("L107").Value = 4125
$ws.Range("M107").Value = -3183
$ws.Range("N107").Value = -7965
